$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = "04/08/2025"
$ws.Range("A19").Style = "Normal"

$ws.Range("B19").Value = "U. De Chile"
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = "Cobresal"
$ws.Range("F19").Value = "L"
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 1
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 1.68
$ws.Range("L19").Value = 0.38
$ws.Range("M19").Value = 22
$ws.Range("N19").Value = 7
$ws.Range("O19").Value = 3
$ws.Range("P19").Value = 4
